# Improve IO Control testbench
# Add a new time-record row (row 7) documenting additional IO Control Unit
# testbench work on 25.11.2019, 22:30-23:00.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Re-use the existing "25.11.2019" value already used in A5/A6 for the date.
$lastDate = $ws.Range("A6").Value()
$ws.Range("A7").Value = $lastDate

# Start / end time of the new entry (22:30 -> 23:00), matching the time
# formatting used by the row above (B6:C6).
$ws.Range("B7").Value = 0.9375
$ws.Range("C7").Value = 0.95833333333333337
$ws.Range("B7:C7").NumberFormat = $ws.Range("B6:C6").NumberFormat

# Duration formula, consistent with the rest of column D.
$ws.Range("D7").Formula = "=C7-B7"
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat

# Category / description for the new entry.
$ws.Range("E7").Value = "IO Control Unit"
$ws.Range("F7").Value = "Testbench"

# Move the active selection past the newly added row, like Excel does
# after entering data in the last row.
$ws.Range("A8").Select()
